$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the "Date" property value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-12T09:15:29+00:00"

# --- Elements sheet: append a new mapping column (AL) ---
$elements = $wb.Worksheets.Item("Elements")

# Duplicate the formatting of the last existing column (AK) into the new
# column (AL) so the header keeps the bold/filled style and the data rows
# keep the bordered body style.
$elements.Range("AK1").Copy()
$elements.Range("AL1").PasteSpecial(-4122) | Out-Null

$elements.Range("AK2:AK6").Copy()
$elements.Range("AL2:AL6").PasteSpecial(-4122) | Out-Null

$elements.Application.CutCopyMode = $false

# Header text for the new mapping column
$elements.Range("AL1").Value = "Mapping: Spécification métier vers l'extension ROR HealthcareServiceNoConsentHabilitation"

# Data rows: only the last row (the new draft mapping) gets a value, the
# rest stay empty just like the other mapping columns for this element set.
$elements.Range("AL2").Value = ""
$elements.Range("AL3").Value = ""
$elements.Range("AL4").Value = ""
$elements.Range("AL5").Value = ""
$elements.Range("AL6").Value = "habilitationAuxSoinsSansConsentement"

# Give the new column a sensible width matching the other wide mapping columns
# (closest achievable value to 97.33203125 via the ColumnWidth property)
$elements.Columns.Item(38).ColumnWidth = 96.5
